$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

Set-TextValue 'D2' '59.859.06'
Set-TextValue 'E2' '  +1.52%  '
Set-TextValue 'D3' '2.636.40'
Set-TextValue 'E4' '  +0.07%  '
Set-TextValue 'D5' '537.84'
Set-TextValue 'E5' '  +1.74%  '
Set-TextValue 'E6' '  +3.87%  '
Set-TextValue 'E7' '  +0.04%  '
Set-TextValue 'D8' '0.568'
Set-TextValue 'E8' '  +0.58%  '
Set-TextValue 'D9' '6.56'
Set-TextValue 'E9' '  +2.45%  '
Set-TextValue 'E10' '  +1.90%  '
Set-TextValue 'E11' '  +1.52%  '
Set-TextValue 'E12' '  -1.71%  '
Set-TextValue 'D13' '3.099.51'
Set-TextValue 'E13' '  +1.97%  '
Set-TextValue 'D14' '59.763.58'
Set-TextValue 'E14' '  +1.46%  '
Set-TextValue 'D15' '20.97'
Set-TextValue 'E15' '  +2.50%  '
Set-TextValue 'D16' '2.619.21'
Set-TextValue 'E16' '  +1.97%  '
Set-TextValue 'E17' '  +1.57%  '
Set-TextValue 'D18' '343.41'
Set-TextValue 'E18' '  -0.12%  '
Set-TextValue 'E19' '  +2.64%  '
Set-TextValue 'D20' '10.23'
Set-TextValue 'E20' '  +1.62%  '
Set-TextValue 'D21' '6.42'
Set-TextValue 'E21' '  -0.28%  '
Set-TextValue 'E22' '  +0.01%  '
Set-TextValue 'D23' '67.57'
Set-TextValue 'E23' '  +1.21%  '
Set-TextValue 'D24' '0.412'
Set-TextValue 'E24' '  +2.01%  '
Set-TextValue 'E25' '  -0.51%  '
Set-TextValue 'E26' '  +0.01%  '
Set-TextValue 'D27' '7.26'
Set-TextValue 'E27' '  +2.93%  '
Set-TextValue 'D28' '0.0₃0752'
Set-TextValue 'E28' '  +4.92%  '
Set-TextValue 'E30' '  +4.18%  '
Set-TextValue 'D31' '5.87'
Set-TextValue 'E31' '  -0.12%  '
Set-TextValue 'D32' '19.00'
Set-TextValue 'E32' '  +1.63%  '
Set-TextValue 'D33' '150.80'
Set-TextValue 'E33' '  +1.43%  '
Set-TextValue 'D34' '4.02'
Set-TextValue 'E34' '  +1.98%  '
Set-TextValue 'D35' '1.14'
Set-TextValue 'E35' '  +2.08%  '
Set-TextValue 'B36' 'Fetch.AI'
Set-TextValue 'C36' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D36' '0.840'
Set-TextValue 'E36' '  +2.23%  '
Set-TextValue 'B37' 'Stacks'
Set-TextValue 'C37' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue 'D37' '1.46'
Set-TextValue 'E37' '  -0.86%  '
Set-TextValue 'E38' '  +1.89%  '
Set-TextValue 'D39' '290.14'
Set-TextValue 'E39' '  +8.29%  '
Set-TextValue 'D40' '3.58'
Set-TextValue 'E40' '  +1.86%  '
Set-TextValue 'E41' '  +0.13%  '
Set-TextValue 'D42' '0.604'
Set-TextValue 'E42' '  +1.07%  '
Set-TextValue 'D43' '10.73'
Set-TextValue 'E43' '  -0.16%  '
Set-TextValue 'D44' '0.0952'
Set-TextValue 'E44' '  -0.03%  '
Set-TextValue 'D45' '0.0533'
Set-TextValue 'E45' '  +3.96%  '
Set-TextValue 'D46' '1.969.19'
Set-TextValue 'E46' '  +0.56%  '
Set-TextValue 'E47' '  +1.97%  '
Set-TextValue 'D48' '18.51'
Set-TextValue 'E48' '  +2.02%  '
Set-TextValue 'E49' '  +3.04%  '
Set-TextValue 'D50' '110.82'
Set-TextValue 'E50' '  -0.92%  '
Set-TextValue 'D51' '4.74'
Set-TextValue 'E51' '  +0.00%  '
